# update: sync, split_day and update payment 20241223 - 20250210
#
# - shuttlecock_buy: append purchase #19 (Yonex AS-10, bought 2025-01-01)
# - Sheet1 (new tab, after shuttlecock_buy): purchase #19's "split day"
#   re-buy entry (Yonex AS-10, bought 2025-01-13)

$wb = $excel.ActiveWorkbook

# --- add the new "Sheet1" tab, placed right after "shuttlecock_buy" ---
$ws2 = $wb.Worksheets.Add()
$ws2.Name = "Sheet1"
$ws2.Move($null, $wb.Worksheets.Item("shuttlecock_buy"))

# Re-fetch the original sheet: Worksheets collection positions shifted
# once the new sheet was inserted/moved, so grab it fresh by name.
$ws1 = $wb.Worksheets.Item("shuttlecock_buy")

# --- shuttlecock_buy: new row 20 (index 19) ---
$ws1.Range("A20").Value = 19
$ws1.Range("B20").Value = 45658
$ws1.Range("B20").NumberFormat = "d-mmm-yy"
$ws1.Range("C20").Value = "Yonex AS-10"
$ws1.Range("D20").Value = 810
$ws1.Range("E20").Value = 10
$ws1.Range("F20").Value = 0
$ws1.Range("G20").Formula = "=D20*E20+F20"
$ws1.Range("H20").Formula = "=E20*12"
$ws1.Range("I20").Formula = "=ROUNDUP(G20/H20,0)"

# --- Sheet1: row 2 (index 19, split-day re-buy) ---
$ws2 = $wb.Worksheets.Item("Sheet1")
$ws2.Range("A2").Value = 19
$ws2.Range("B2").Value = 45670
$ws2.Range("B2").NumberFormat = "d-mmm-yy"
$ws2.Range("C2").Value = "Yonex AS-10"
$ws2.Range("D2").Value = 825
$ws2.Range("E2").Value = 10
$ws2.Range("F2").Value = 0
$ws2.Range("G2").Formula = "=D2*E2+F2"
$ws2.Range("H2").Formula = "=E2*12"
$ws2.Range("I2").Formula = "=ROUNDUP(G2/H2,0)"

# --- match the final selection recorded on each sheet ---
$ws2.Range("D13").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("N10").Select() | Out-Null
